$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 30: continuation entry after row 29 (date 41205 -> 41206, i.e. 2012-10-23 -> 2012-10-24)
$ws.Range("A30").Value = [DateTime]::FromOADate(41206)
$ws.Range("A30").NumberFormat = 'ddd\ dd/mm/yyyy'

$ws.Range("B30").Value = 2.75

$ws.Range("D30").Value = "Manual continued"

$ws.Range("D30").Select()
